$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting (styles) of an existing, similarly-structured
# row (row 21 uses the same per-column style pattern as the new row) onto
# the new row 27 before filling in the real values.
$ws.Range("A21:G21").Copy($ws.Range("A27:G27"))

# Populate the new review row.
$ws.Range("A27").Value = "passive.income.nadi.myfirstdrawermenuproject"
$ws.Range("B27").Value = "passive income"
$ws.Range("C27").Value = "itaisenior@gmail.com"
$ws.Range("D27").Value = "vikicrestina@gmail.com"
$ws.Range("E27").Value = "27/5/2019 15:59"
$ws.Range("F27").Value = "this info is pricless, especially the secret"
$ws.Range("G27").Value = "no"

# Wire up the two mailto hyperlinks for the new row.
$ws.Hyperlinks.Add($ws.Range("C27"), "mailto:itaisenior@gmail.com", $null, $null, "itaisenior@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D27"), "mailto:vikicrestina@gmail.com", $null, $null, "vikicrestina@gmail.com") | Out-Null

# Hyperlinks.Add() re-styles the target cells with the built-in Hyperlink
# style; restore the original column formatting that the rest of the
# sheet uses (copied above) now that the relationships are established.
$ws.Range("C21:D21").Copy($ws.Range("C27:D27"))

# Move the active selection the way the source workbook left it.
$ws.Range("G28").Select() | Out-Null
